# Rebuild the "Repartition Des Taches" planning table on Feuil1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Tâches / Fonctions"
$ws.Range("B1").Value = "Responsable"
$ws.Range("C1").Value = "Début"
$ws.Range("D1").Value = "Fin"
$ws.Range("E1").Value = "Statut"

$ws.Range("A2").Value = "Création Des Dossiers"
$ws.Range("B2").Value = "Yoan / Hugo"
$ws.Range("C2").Value = 44907
$ws.Range("C2").NumberFormat = "d-mmm"
$ws.Range("D2").Value = 44913
$ws.Range("D2").NumberFormat = "d-mmm"
$ws.Range("E2").Interior.Color = 5287936

$ws.Range("A3").Value = "Création Du Template"
$ws.Range("B3").Value = "Yoan / Hugo"
$ws.Range("C3").Value = 44907
$ws.Range("C3").NumberFormat = "d-mmm"
$ws.Range("D3").Value = 44907
$ws.Range("D3").NumberFormat = "d-mmm"
$ws.Range("E3").Interior.Color = 5287936

$ws.Range("A4").Value = "Modification Hangman.go"
$ws.Range("B4").Value = "Yoan / Hugo"
$ws.Range("C4").Value = 44907
$ws.Range("C4").NumberFormat = "d-mmm"
$ws.Range("D4").Value = 44913
$ws.Range("D4").NumberFormat = "d-mmm"
$ws.Range("E4").Interior.Color = 5287936

$ws.Range("A5").Value = "Création du CSS"
$ws.Range("B5").Value = "Hugo"
$ws.Range("C5").Value = 44908
$ws.Range("C5").NumberFormat = "d-mmm"
$ws.Range("D5").Value = 44913
$ws.Range("D5").NumberFormat = "d-mmm"
$ws.Range("E5").Interior.Color = 5287936

$ws.Range("A6").Value = "Création Fonctions go"
$ws.Range("B6").Value = "Yoan"
$ws.Range("C6").Value = 44908
$ws.Range("C6").NumberFormat = "d-mmm"
$ws.Range("D6").Value = 44913
$ws.Range("D6").NumberFormat = "d-mmm"
$ws.Range("E6").Interior.Color = 5287936

$ws.Range("A7").Value = "Création HTML"
$ws.Range("B7").Value = "Yoan / Hugo"
$ws.Range("C7").Value = 44908
$ws.Range("C7").NumberFormat = "d-mmm"
$ws.Range("D7").Value = 44913
$ws.Range("D7").NumberFormat = "d-mmm"
$ws.Range("E7").Interior.Color = 5287936

$ws.Range("A8").Value = "Demande Pseudo Utilisateur"
$ws.Range("B8").Value = "Yoan "
$ws.Range("C8").Value = 44907
$ws.Range("C8").NumberFormat = "d-mmm"
$ws.Range("D8").Value = 44907
$ws.Range("D8").NumberFormat = "d-mmm"
$ws.Range("E8").Interior.Color = 5287936

$ws.Range("A9").Value = "Choix Niveau"
$ws.Range("B9").Value = "Hugo"
$ws.Range("C9").Value = 44908
$ws.Range("C9").NumberFormat = "d-mmm"
$ws.Range("D9").Value = 44908
$ws.Range("D9").NumberFormat = "d-mmm"
$ws.Range("E9").Interior.Color = 5287936

$ws.Range("A10").Value = "Jeu Hangman Web"
$ws.Range("B10").Value = "Yoan / Hugo"
$ws.Range("C10").Value = 44908
$ws.Range("C10").NumberFormat = "d-mmm"
$ws.Range("D10").Value = 44913
$ws.Range("D10").NumberFormat = "d-mmm"
$ws.Range("E10").Interior.Color = 5287936

$ws.Range("A11").Value = "Dessin du pendu"
$ws.Range("B11").Value = "Hugo"
$ws.Range("C11").Value = 44912
$ws.Range("C11").NumberFormat = "d-mmm"
$ws.Range("D11").Value = 44913
$ws.Range("D11").NumberFormat = "d-mmm"
$ws.Range("E11").Interior.Color = 5287936

$ws.Range("A12").Value = "Condition Victoire / defaite"
$ws.Range("B12").Value = "Yoan"
$ws.Range("C12").Value = 44911
$ws.Range("C12").NumberFormat = "d-mmm"
$ws.Range("D12").Value = 44913
$ws.Range("D12").NumberFormat = "d-mmm"
$ws.Range("E12").Interior.Color = 5287936

$ws.Range("A13").Value = "Readme"
$ws.Range("B13").Value = "Yoan / Hugo"
$ws.Range("C13").Value = 44913
$ws.Range("C13").NumberFormat = "d-mmm"
$ws.Range("D13").Value = 44913
$ws.Range("D13").NumberFormat = "d-mmm"
$ws.Range("E13").Interior.Color = 5287936


# Column widths (best achievable approximation given engine's width quantization)
$ws.Columns.Item(1).ColumnWidth = 24.666666666666668
$ws.Columns.Item(2).ColumnWidth = 11.5
$ws.Columns.Item(3).ColumnWidth = 8.333333333333334
$ws.Columns.Item(4).ColumnWidth = 17.833333333333332

# Selection & zoom state to match the saved view
$excel.ActiveWindow.Zoom = 205
$ws.Range("E13").Select()
